$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. New row: méthodes principales de gestion ... -> X for PIGUET Charles (col 3) ---
$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "Ajout des méthodes principales de gestion des épreuves, des organisateurs, des inscriptions et des événements dans le contrôleur"
$row1.Cells.Item(3).Range.Text = "X"

# --- 2. New row: Fonctionnalité d'upload du CSV -> X for DALICHAMPT Thibaut (col 2) ---
$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Fonctionnalité d’upload du fichier CSV du classement d’une épreuve"
$row2.Cells.Item(2).Range.Text = "X"

# --- 3. New row: Agencement de la grille et du design SCSS -> X for CLEMENT Dylan (col 5) ---
$row3 = $t.Rows.Add()
$row3.Cells.Item(1).Range.Text = "Agencement de la grille et du design SCSS"
$row3.Cells.Item(5).Range.Text = "X"
